# Append: 2026-01-14 12:41 JST
# This script updates the "ランサーズ" sheet (sheet1) of the workbook:
#  - refreshes the timestamp in column A for every existing data row
#  - inserts 4 new job-listing rows at their correct positions
#  - rebuilds the hyperlinks for column F in the new row order

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks up front; they will be re-added at the end
# once every row is in its final location (Rows.Insert does not shift the
# <hyperlinks> relationships automatically).
$ws.Hyperlinks.Delete()

# --- Insert 4 blank rows at the positions where brand-new listings land ---
# Doing this top-to-bottom, each insertion shifts everything below it down
# by one row, so later insert positions already account for earlier ones.
$ws.Rows.Item(3).Insert()   # new row 3 (was nothing yet)
$ws.Rows.Item(4).Insert()   # new row 4
$ws.Rows.Item(7).Insert()   # new row 7
$ws.Rows.Item(9).Insert()   # new row 9

$timestamp = "2026-01-14 12:41:23"

# --- Row 2: unchanged listing, only the fetch timestamp changes ---
$ws.Cells.Item(2, 1).Value = $timestamp

# --- Row 3 (new): AI operations designer listing ---
$ws.Cells.Item(3, 1).Value = $timestamp
$ws.Cells.Item(3, 2).Value = "AIオペレーションデザイナーを募集します(経営直結/会議→意思決定変換)"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5471032"
$ws.Cells.Item(3, 6).Style = "Hyperlink"
$ws.Cells.Item(3, 7).Value = 310
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai"

# --- Row 4 (new): Google business profile bulk tool listing ---
$ws.Cells.Item(4, 1).Value = $timestamp
$ws.Cells.Item(4, 2).Value = "googleビジネスプロフィール一括ツール"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5470814"
$ws.Cells.Item(4, 6).Style = "Hyperlink"
$ws.Cells.Item(4, 7).Value = 73
$ws.Cells.Item(4, 8).Value = "◆ツール"

# --- Row 5: existing listing (was row 3), timestamp refresh only ---
$ws.Cells.Item(5, 1).Value = $timestamp

# --- Row 6: existing listing (was row 4), timestamp refresh only ---
$ws.Cells.Item(6, 1).Value = $timestamp

# --- Row 7 (new): Online web reservation system listing ---
$ws.Cells.Item(7, 1).Value = $timestamp
$ws.Cells.Item(7, 2).Value = "オンラインWEB予約システム構築"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5470812"
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 7).Value = 33

# --- Row 8: existing listing (was row 5), timestamp refresh only ---
$ws.Cells.Item(8, 1).Value = $timestamp

# --- Row 9 (new): Moving enclosure control harness listing ---
$ws.Cells.Item(9, 1).Value = $timestamp
$ws.Cells.Item(9, 2).Value = "移動型筐体の制御ハーネス製作(Arduino/電飾/音声/電源)"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5471022"
$ws.Cells.Item(9, 6).Style = "Hyperlink"
$ws.Cells.Item(9, 7).Value = 18

# --- Row 10: existing listing (was row 6), timestamp refresh only ---
$ws.Cells.Item(10, 1).Value = $timestamp

# --- Row 11: existing listing (was row 7), timestamp refresh only ---
$ws.Cells.Item(11, 1).Value = $timestamp

# --- Rebuild hyperlinks for F2:F11 in final row order ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5470737")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5471032")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5470814")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5470623")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5470812")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5470403")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5471022")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5470150")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5470726")

# Re-apply the hyperlink visual style to every link cell (Hyperlinks.Add can
# reset formatting) and make sure F2's style (already "Hyperlink" originally)
# stays consistent too.
$ws.Range("F2:F11").Style = "Hyperlink"

[void]$ws.Range("A1").Select()
